# Customer invoice: decrease table header font by 1
#
# The first table's header row (row 1) currently has every run/paragraph
# mark set to w:sz=24 (12pt) / w:szCs=24 (12pt). The commit shrinks that
# header row's font by one point: w:sz becomes 20 (10pt) and w:szCs
# becomes 22 (11pt).

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)
$headerRow = $table.Rows.Item(1)
$rng = $headerRow.Range

$rng.Font.Size = 10
$rng.Font.SizeBi = 11
